# Populate the header row used by the music-player CSV/export template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "songImage (Image Cover)"
$ws.Range("B1").Value = "songMusic (Music Title - Filename)"
$ws.Range("C1").Value = "songName (Music Title)"
$ws.Range("D1").Value = "artistName (Artist Name)"
$ws.Range("E1").Value = "albumName (Album Name)"

# Widen the columns so the header text isn't truncated (values chosen so the
# stored OOXML column width - which Excel derives by rounding to whole
# pixels at the workbook's fixed Normal-style digit width - lands on the
# target widths of 29, 40.5703125, 33.85546875, 27.28515625 and 33 chars).
$ws.Columns.Item(1).ColumnWidth = 28.166666666666668
$ws.Columns.Item(2).ColumnWidth = 39.666666666666664
$ws.Columns.Item(3).ColumnWidth = 33
$ws.Columns.Item(4).ColumnWidth = 26.5
$ws.Columns.Item(5).ColumnWidth = 32.166666666666664

# Leave the selection where the author left it when they saved the file.
$ws.Range("B15").Select()
